$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$regexIntro = "The syntax impl is most similar to perl among all differnet regex util like grep, awk and etc. Three major class is used in java.util.regex package:`n1. Pattern: represent the compiled pattern, construct by static factory method Pattern.compile(String regex)`n2. Matcher: represent the one who has the knowlege to match string with pattern. its construct by objPattern.matcher(String str_to_search)`n3. PatternSyntaxException"

$regexSyntaxTitle = "Regex: Syntax in java"
$regexExampleTitle = "Regex: Example"

$regexExample = "Pattern pattern = Pattern.compile(`"my.*[a-z]+`");`nMatcher matcher = pattern.matcher(`"A string to be searched`");`nboolean found = false;`n while (matcher.find()) {`n                console.format(`"I found the text`" +`n                    `" \`"%s\`" starting at `" +`n                    `"index %d and ending at index %d.%n`",`n                    matcher.group(),`n                    matcher.start(),`n                    matcher.end());`n                found = true;`n }"

$regexSyntaxManual = "Below topics are covered in the manual:`n1.  regex of Character ([a-z])`n2.  regex of Predefined Character (\d, \D)`n3. regex of qualifier (X?, X*, X+)`n4. Capturing Group ((A)*(B)(C)+)`n5. regex of boundary (^, `$)`n6. manual of Pattern`n7. manual of Matcher`n8. manual of XXXException`n9. how to match unicode (\uxxxx)"

# Row 114: RTFM / Regex / intro text
$ws.Range("A114").Value = "RTFM"
$ws.Range("B114").Value = "Regex"
$ws.Range("C114").Value = $regexIntro

# Row 115: RTFM / Regex: Syntax in java / (syllabus further down)
$ws.Range("A115").Value = "RTFM"
$ws.Range("B115").Value = $regexSyntaxTitle

# Row 116: RTFM / Regex: Example / code example
$ws.Range("A116").Value = "RTFM"
$ws.Range("B116").Value = $regexExampleTitle
$ws.Range("C116").Value = $regexExample

# Fill in C115 last so the new shared-string creation order matches
$ws.Range("C115").Value = $regexSyntaxManual

$ws.Rows.Item(114).RowHeight = 32.25
$ws.Rows.Item(115).RowHeight = 32.25
$ws.Rows.Item(116).RowHeight = 32.25

$ws.Range("A116").Select()
